# Generate Report for handback
# Updates the localization-status workbook: marks zh-cn/de-de rows 2-3 as
# "Handed back: in sync with en-US", fills in the Latest Target File /
# Latest Handback File columns (E/F) with hyperlinks, and stamps the
# Latest Handback DateTime (column G) with the handback timestamp.

$wb = $excel.ActiveWorkbook

$mdDisplay   = "4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.md"
$mdUrl       = "https://github.com/OpenLocalizationTest/oltest/blob/17cadf0d84543e7bc3e62da3e9ac5e1cb1597fbb/e2e/4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.md"

$statusText  = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$zhXlfDisplay = "4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.d76442c13d83cd579ed92490c6b7780c957ec87e.zh-cn.xlf"
$zhXlfUrl     = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/daa45f64cde85a6f199e279e900d5e46f565234b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.d76442c13d83cd579ed92490c6b7780c957ec87e.zh-cn.xlf"
$zhHandback   = "2016-01-28 09:38:09"

$wsZh.Range("B2").Value = $statusText
$wsZh.Range("B3").Value = $statusText

$wsZh.Hyperlinks.Add($wsZh.Range("E2"), $mdUrl, "", "", $mdDisplay)
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $zhXlfUrl, "", "", $zhXlfDisplay)
$wsZh.Hyperlinks.Add($wsZh.Range("E3"), $mdUrl, "", "", $mdDisplay)
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), $zhXlfUrl, "", "", $zhXlfDisplay)

$wsZh.Range("G2").Value = $zhHandback
$wsZh.Range("G3").Value = $zhHandback

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$deXlfDisplay = "4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.d76442c13d83cd579ed92490c6b7780c957ec87e.de-de.xlf"
$deXlfUrl     = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/426b6feca0c0f35facdbef193cd977f5b1d3718b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/4d0d0e8a-fe96-4ce6-bad6-7e1f3d85b319.d76442c13d83cd579ed92490c6b7780c957ec87e.de-de.xlf"
$deHandback   = "2016-01-28 09:38:30"

$wsDe.Range("B2").Value = $statusText
$wsDe.Range("B3").Value = $statusText

$wsDe.Hyperlinks.Add($wsDe.Range("E2"), $mdUrl, "", "", $mdDisplay)
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $deXlfUrl, "", "", $deXlfDisplay)
$wsDe.Hyperlinks.Add($wsDe.Range("E3"), $mdUrl, "", "", $mdDisplay)
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), $deXlfUrl, "", "", $deXlfDisplay)

$wsDe.Range("G2").Value = $deHandback
$wsDe.Range("G3").Value = $deHandback
